# Update countries & provincias Spain
# - Reorders a few country rows (name swaps with updated stats) and
#   refreshes case counts for several countries in the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Estados Unidos (row 4): refresh stats ---
$ws.Range("B4").Value = 1184711
$ws.Range("C4").Value = 23937
$ws.Range("E4").Value = 938003
$ws.Range("G4").Value = 1045
$ws.Range("H4").Value = 68489

# --- Brasil (row 12): refresh stats ---
$ws.Range("B12").Value = 101147
$ws.Range("C12").Value = 4588
$ws.Range("E12").Value = 53185
$ws.Range("G12").Value = 275
$ws.Range("H12").Value = 7025

# --- Canada (row 15): refresh stats ---
$ws.Range("B15").Value = 59378
$ws.Range("C15").Value = 2664
$ws.Range("D15").Value = 24729
$ws.Range("E15").Value = 30968
$ws.Range("G15").Value = 115
$ws.Range("H15").Value = 3681

# --- Israel (row 31): refresh stats ---
$ws.Range("B31").Value = 16208
$ws.Range("C31").Value = 23
$ws.Range("D31").Value = 9749
$ws.Range("E31").Value = 6227
$ws.Range("G31").Value = 3
$ws.Range("H31").Value = 232

# --- Burkina Faso moves above Uruguay (row 107/108), both get fresh data ---
$ws.Range("A107").Value = "Burkina Faso"
$ws.Range("B107").Value = 662
$ws.Range("C107").Value = 10
$ws.Range("D107").Value = 540
$ws.Range("E107").Value = 77
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 45

$ws.Range("A108").Value = "Uruguay"
$ws.Range("B108").Value = 652
$ws.Range("C108").Value = 0
$ws.Range("D108").Value = 440
$ws.Range("E108").Value = 195
$ws.Range("F108").Value = 10
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 17

# --- Ruanda moves above Guinea-Bisau (row 132/133), both get fresh data ---
$ws.Range("A132").Value = "Ruanda"
$ws.Range("B132").Value = 259
$ws.Range("C132").Value = 4
$ws.Range("D132").Value = 124
$ws.Range("E132").Value = 135
$ws.Range("H132").Value = 0

$ws.Range("A133").Value = "Guinea-Bisau"
$ws.Range("B133").Value = 257
$ws.Range("D133").Value = 19
$ws.Range("E133").Value = 237
$ws.Range("H133").Value = 1

# --- Sierra Leona (row 137): refresh stats ---
$ws.Range("B137").Value = 166
$ws.Range("C137").Value = 11
$ws.Range("E137").Value = 129

# --- Liberia moves above Birmania (row 138/139), both get fresh data ---
$ws.Range("A138").Value = "Liberia"
$ws.Range("B138").Value = 158
$ws.Range("D138").Value = 58
$ws.Range("E138").Value = 82
$ws.Range("H138").Value = 18

$ws.Range("A139").Value = "Birmania"
$ws.Range("B139").Value = 155
$ws.Range("C139").Value = 4
$ws.Range("D139").Value = 43
$ws.Range("E139").Value = 106
$ws.Range("H139").Value = 6

# --- Namibia moves above San Vicente y las Granadinas (row 194/195) ---
$ws.Range("A194").Value = "Namibia"
$ws.Range("A195").Value = "San Vicente y las Granadinas"

Write-Host "Applied country/provincia updates"
